$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Tanggal Selesai" for Preprocessing Data (row 3): 2021-10-27 -> 2021-10-28
$ws.Range("C3").Value = [datetime]"2021-10-28"

# Update "Tanggal Mulai" / "Tanggal Selesai" for Modelling CNN (row 4)
$ws.Range("B4").Value = [datetime]"2021-10-28"
$ws.Range("C4").Value = [datetime]"2021-10-31"

# Move the active selection to C5, as in the saved file
$ws.Range("C5").Select()
